$wb = $excel.ActiveWorkbook

$rowsToUpdate = @(7, 8, 9, 11, 12, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) timestamp update ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rowsToUpdate) {
    $wsOverview.Range("G$r").Value = "2016-09-02 06:26:35"
}

# --- zh-cn sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-02 06:26:29"
}

# --- de-de sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-02 06:26:35"
}
